# Fill in the template with real example data, add a hyperlink for the
# email address, format the phone number as text, widen column B and
# move the active selection down to B6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the placeholder instructions with actual sample values.
$ws.Range("B1").Value = "John"
$ws.Range("B2").Value = "+48521439678"
$ws.Range("B3").Value = "john@test.com"
$ws.Range("B4").Value = "name of your team"
$ws.Range("B5").Value = "number of players"
$ws.Range("B6").Value = "password from your mzgb account"

# Store phone number and team name as text (quote-prefixed) and add a
# hyperlink on the email cell.
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B4").NumberFormat = "@"

$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:john@test.com")

# Widen column B slightly to fit the new content.
$ws.Columns.Item(2).ColumnWidth = 35.5703125

# Move the selected cell to B6.
$ws.Range("B6").Select()

$wb.Save()
